# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values, columns A..AG
$ws.Range("A2").Value = 45977
$ws.Range("B2").Value = 35.32
$ws.Range("C2").Value = 29.46
$ws.Range("D2").Value = 25.77
$ws.Range("E2").Value = 23.58
$ws.Range("F2").Value = 22.95
$ws.Range("G2").Value = 23.08
$ws.Range("H2").Value = 24.09
$ws.Range("I2").Value = 31.08
$ws.Range("J2").Value = 22.97
$ws.Range("K2").Value = 15.83
$ws.Range("L2").Value = 3.03
$ws.Range("M2").Value = 0.25
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0.01
$ws.Range("P2").Value = 0.25
$ws.Range("Q2").Value = 3.64
$ws.Range("R2").Value = 20.07
$ws.Range("S2").Value = 52.24
$ws.Range("T2").Value = 72.61
$ws.Range("U2").Value = 82.44
$ws.Range("V2").Value = 91.59
$ws.Range("W2").Value = 96.13
$ws.Range("X2").Value = 78.33
$ws.Range("Y2").Value = 64.64
$ws.Range("Z2").Value = 34.14
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 82.67
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 93.86
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 77.53
$ws.Range("AG2").Value = "1h-16h"
